# Update the "LOS Galacticos" fantasy roster sheet with the new player
# lineup, positions and teams.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Each row: Player (A), Position (B), Team (C)
$data = @(
    @("Keyonte George",        "PG,SG",    "Utah Jazz"),
    @("Paul George",           "SG,SF,PF", "Philadelphia 76ers"),
    @("James Harden",          "PG,SG",    "LA Clippers"),
    @("Anfernee Simons",       "PG,SG",    "Portland Trail Blazers"),
    @("Jayson Tatum",          "SF,PF",    "Boston Celtics"),
    @("Draymond Green",        "PF,C",     "Golden State Warriors"),
    @("Amen Thompson",         "SG,SF,PF", "Houston Rockets"),
    @("Zion Williamson",       "PF,C",     "New Orleans Pelicans"),
    @("Bobby Portis",          "PF,C",     "Milwaukee Bucks"),
    @("Guerschon Yabusele",    "PF,C",     "Philadelphia 76ers"),
    @("Jaren Jackson Jr.",     "PF,C",     "Memphis Grizzlies"),
    @("Moussa Diabate",        "C",        "Charlotte Hornets"),
    @("Ivica Zubac",           "C",        "LA Clippers"),
    @("Bradley Beal",          "PG,SG,SF", "Phoenix Suns"),
    @("Anthony Edwards",       "SG,SF",    "Minnesota Timberwolves"),
    @("Giannis Antetokounmpo", "PF,C",     "Milwaukee Bucks"),
    @("Fred VanVleet",         "PG",       "Houston Rockets"),
    @("RJ Barrett",            "SG,SF,PF", "Toronto Raptors")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
